$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2005
$ws1.Range("F7").Value = 3354
$ws1.Range("F9").Value = 811

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 2005
$ws4.Range("F8").Value = 3354
$ws4.Range("F10").Value = 811
